$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# --- Row 11 (Sno 10): append "order by p.LST_CHG_BY_DTTM desc" to the portal_user_tin query ---
$ws.Range("B11").Value = "select * from OLE.PORTAL_USER p join OLE.portal_user_tin pt `non p.PORTAL_USER_ID=pt.PORTAL_USER_ID`nwhere pt.PROV_TIN_NBR='{`$tin}' and p.STS_CD='A' and p.USERNAME not in( '{`$id}' ) and p.USERNAME not in ('USPROV7729','AUAUTO1563','AUAUTO5896','AUBILL5903','AUPAYE5365') order by p.LST_CHG_BY_DTTM desc`n"
$ws.Rows.Item(11).RowHeight = 86.4

# --- Row 19 (Sno 18): append "order by p.LST_CHG_BY_DTTM desc Fetch FIRST ROW ONLY" to the BS_TIN query ---
$ws.Range("A19").Value = "'18"
$ws.Range("B19").Value = "  `nselect * from OLE.PORTAL_USER p join OLE.PORTAL_USER_BS_TIN bt `non p.PORTAL_USER_ID=bt.PORTAL_USER_ID join OLE.BILLING_SERVICE bs`n on bt.BILLING_SERVICE_ID=bs.BILLING_SERVICE_ID`nwhere  bs.IDENTIFIER_NBR='{`$tin}'  and p.STS_CD='A' and p.USERNAME not in ('USPROV7729','AUAUTO1563','AUAUTO5896','AUBILL5903','AUPAYE5365') order by p.LST_CHG_BY_DTTM desc Fetch FIRST ROW ONLY"
$ws.Rows.Item(19).RowHeight = 100.8

# --- Row 20 (Sno 19, new): payer TIN lookup query ---
$ws.Range("A20").Value = "'19"
$ws.Range("B20").Value = "select * from OLE.PORTAL_USER p join`nOLE.PORTAL_USER_PAYER_TIN py `non p.PORTAL_USER_ID=py.PORTAL_USER_ID `nwhere  py.PAYR_TIN_NBR='{`$tin}' and p.STS_CD='A' and p.USERNAME not in ('USPROV7729','AUAUTO1563','AUAUTO5896','AUBILL5903','AUPAYE5365') order by p.LST_CHG_BY_DTTM desc Fetch FIRST ROW ONLY"
$ws.Rows.Item(20).RowHeight = 86.4

# --- Scroll the view down a bit (top-left cell moves from A11 to A13) ---
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B19").Select()

Write-Host "Edit applied"
